$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12's "Approved/Rejected" result is corrected from "Rejected" to "Approved",
# and the stray "Teststep wrong" reason in the "ReasonToReject" column is cleared.
$ws.Range("I12").Value = "Approved"
$ws.Range("J12").ClearContents()

# Update the selection to match the new active cell.
$ws.Range("J12").Select()
